$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: remove the "Meta description: ..." paragraph (the second
# paragraph in the document, right after the H1 title).
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# Change 2: the final paragraph used to hold the image-generation "Prompt:"
# text in italics. It now becomes two paragraphs:
#   1) a new bold paragraph repeating the page title
#   2) the old paragraph, still italic, but with its text replaced by the
#      meta-description copy we removed above.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

# First, swap the "Prompt: ..." text for the meta-description text while the
# paragraph's own (italic) run formatting stays untouched.
$promptText = "Prompt: Create a fun and engaging feature image for Book of Darkness that captures the thrilling concept of the game. The image should be in a cartoon style and feature a happy Maya warrior wearing glasses, surrounded by fire and magic symbols. Make sure to include the iconic Book of Darkness in the image as well. Overall, the image should showcase the power struggle between the good and evil forces in the game, with the Maya warrior as the hero fighting against the evil magician. The image should be vibrant and eye-catching, with bold colors and dynamic imagery that draws in the player's attention."
$replacementText = "Read our comprehensive review of Book of Darkness slot game. Discover exciting gameplay features and high RTP. Play now for free!"

$lastPara.Range.Find.Execute($promptText, $true, $false, $false, $false, $false, $true, 1, $false, $replacementText, 2)

# Now insert a brand-new bold paragraph right before it, carrying the page
# title text.
$pPrev = $d.Paragraphs.Item($count - 1)
$rng = $pPrev.Range
$rng.Collapse(0)
$titleText = "Play Book of Darkness Free - Review of Gameplay Features and More"
$rng.InsertAfter($titleText + "`r")

$count2 = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($count2 - 1)
$titleRng = $newPara.Range
$titleRng.MoveEnd(1, -1)
$titleRng.Font.Bold = 1
